$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 50, shifting existing rows 50-66 down to 51-67
$ws.Rows.Item(50).Insert()

# Fill the new row 50 with the new data point
$ws.Cells.Item(50, 1).Value = 6
$ws.Cells.Item(50, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(50, 3).Value = "Metropolitana"
$ws.Cells.Item(50, 4).Value = 44627
$ws.Cells.Item(50, 4).NumberFormat = $ws.Cells.Item(51, 4).NumberFormat
$ws.Cells.Item(50, 5).Value = 13
$ws.Cells.Item(50, 6).Value = 100114007
$ws.Cells.Item(50, 7).Value = "Jengibre"
$ws.Cells.Item(50, 8).Value = "Sin especificar"
$ws.Cells.Item(50, 9).Value = "Primera"
$ws.Cells.Item(50, 10).Value = 180
$ws.Cells.Item(50, 11).Value = 14000
$ws.Cells.Item(50, 12).Value = 15000
$ws.Cells.Item(50, 13).Value = 14444
$ws.Cells.Item(50, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(50, 15).Value = "Perú"
$ws.Cells.Item(50, 16).Value = 1111
$ws.Cells.Item(50, 17).Value = 13
$ws.Cells.Item(50, 18).Value = "Hortaliza"
